$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.472.07"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "2.163.90"
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'228.07"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("D7").Value = "'63.96"
$ws.Range("E7").Value = "  +2.90%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +2.46%  "

$ws.Range("D10").Value = "'0.0860"
$ws.Range("E10").Value = "  +2.23%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "'16.15"
$ws.Range("E12").Value = "  +2.03%  "

$ws.Range("D13").Value = "2.484.64"
$ws.Range("E13").Value = "  +3.15%  "

$ws.Range("D14").Value = "'22.15"
$ws.Range("E14").Value = "  +0.08%  "

$ws.Range("D15").Value = "'0.814"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D17").Value = "2.136.06"
$ws.Range("E17").Value = "  +1.98%  "

$ws.Range("D18").Value = "39.473.45"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").Value = "'71.92"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "'6.13"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  +1.78%  "

$ws.Range("D22").Value = "'229.92"
$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("D24").Value = "'2.34"
$ws.Range("E24").Value = "  -0.57%  "

$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").Value = "'172.59"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").Value = "'9.52"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "'19.88"
$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  +4.95%  "

$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").Value = "'7.17"
$ws.Range("E34").Value = "  +8.89%  "

$ws.Range("D35").Value = "'4.72"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").Value = "'2.44"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("D38").Value = "'3.56"
$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("E41").Value = "  +0.78%  "

$ws.Range("D42").Value = "'17.67"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("D43").Value = "1.530.86"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").Value = "'1.19"
$ws.Range("E44").Value = "  +3.82%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0933"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.31"
$ws.Range("E46").Value = "  +4.81%  "

$ws.Range("D47").Value = "'2.83"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("E48").Value = "  +5.38%  "

$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").Value = "'9.19"
$ws.Range("E50").Value = "  +24.57%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.367.76"
$ws.Range("E51").Value = "  +3.31%  "
